$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New course rows (2-5) appended below the header row on the "courses" sheet.
# Columns: A vetCode | B cricosCode | C department | D name | E duration
#          H durationDetail | I tuition | J tuitionDetail | M location
# ---------------------------------------------------------------------------

# -- Row 2: DIPLOMA OF CIVIL CONSTRUCTION DESIGN ---------------------------
$ws.Range("A2").Value = "AHC30716"
$ws.Range("B2").Value = "110597F"
$ws.Range("D2").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("E2").Value = 52
$ws.Range("H2").Value = "44 wks Tuition + 8 wks Break"
$ws.Range("M2").Value = "TAS"
$ws.Range("I2").Value = 16200
$ws.Range("J2").Value = "16,000 tuition fee + 200 handling fee"

# -- Row 3: ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN ------------------
$ws.Range("D3").Value = "ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("B3").Value = "111826A"
$ws.Range("A3").Value = "RII60520"
$ws.Range("E3").Value = 104
$ws.Range("H3").Value = "88 wks tuition + 16 wks break"
$ws.Range("I3").Value = 25200
$ws.Range("J3").Value = "25,000 tuition fee + 200 handling fee"
$ws.Range("M3").Value = "TAS"

# -- Row 4: ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY (TELECOMMUNICATIONS NETWORK ENGINEERING)
$ws.Range("D4").Value = "ADVANCED DIPLOMA OF INFORMATION TECHNOLOGY `n(TELECOMMUNICATIONS NETWORK ENGINEERING) "
$ws.Range("A4").Value = "ICT60220"
$ws.Range("B4").Value = "111825B"
$ws.Range("E4").Value = 104
$ws.Range("I4").Value = 13200
$ws.Range("J4").Value = "13,000 tuition fee + 200 handling fee"
$ws.Range("H4").Value = "88 wks tuition + 16 wks break"
$ws.Range("M4").Value = "TAS"

# -- Row 5: PACKAGES (DIPLOMA + ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN)
$ws.Range("C5").Value = "PACKAGES"
$ws.Range("D5").Value = "DIPLOMA OF CIVIL CONSTRUCTION DESIGN + ADVANCED DIPLOMA OF CIVIL CONSTRUCTION DESIGN"
$ws.Range("A5").Value = "RII50520/RII60520"
$ws.Range("B5").Value = "111827M/111826A"
$ws.Range("E5").Value = 104
$ws.Range("I5").Value = 27200
$ws.Range("J5").Value = "27,000 tuition fee + 200 handling fee"
$ws.Range("H5").Value = "88 wks tuition + 16 wks break"
$ws.Range("M5").Value = "TAS"

# -- department column, filled in last across all four rows ----------------
$ws.Range("C2").Value = "CIVIL CONSTRUCTION DESIGN"
$ws.Range("C4").Value = "INFORMATION TECHNOLOGY"
$ws.Range("C3").Value = "CIVIL CONSTRUCTION DESIGN"

# ---------------------------------------------------------------------------
# Formatting: wrap long detail text, apply #,##0 number format to tuition fees
# ---------------------------------------------------------------------------
$ws.Range("H2:H5").WrapText = $true
$ws.Range("D4").WrapText = $true
$ws.Range("A5").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Range("D5").WrapText = $true

$ws.Range("I2:I5").NumberFormat = "#,##0"
$ws.Range("J2:J5").NumberFormat = "#,##0"
$ws.Range("J2:J5").WrapText = $true

# Row heights for the newly-added rows
$ws.Range("A2:R2").RowHeight = 45
$ws.Range("A3:R3").RowHeight = 45
$ws.Range("A4:R4").RowHeight = 45
$ws.Range("A5:R5").RowHeight = 45

# Restore the selection left behind by the editing session
[void]$ws.Range("D18").Select()
